$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared template applied to every new row (14-20): column letter -> value.
# These mirror the constant columns already used by row 13 (inverter row),
# except AG/AH ("meterInterface"/"meterKind") which are "NONE" for the new
# pv-panel rows instead of the firmware version used by the inverter row.
$rowTemplate = [ordered]@{
    'C' = 0
    'D' = 'NONE'
    'E' = 'NONE'
    'F' = 'CMNET'
    'G' = 223
    'H' = 5
    'I' = 5
    'J' = 5
    'K' = 0
    'L' = 0
    'M' = 0
    'N' = 0
    'O' = 0
    'P' = 1
    'Q' = '10F872226797'
    'R' = 0
    'S' = 0
    'T' = 0
    'U' = 0
    'V' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 72
    'AB' = 77
    'AC' = 114
    'AD' = 34
    'AE' = 103
    'AF' = 151
    'AG' = 'NONE'
    'AH' = 'NONE'
    'AI' = 1
    'AJ' = 3600
    'AK' = 'dataeu.hoymiles.com'
    'AL' = 1
    'AM' = 10081
    'AN' = 0
    'AO' = 0
    'AP' = 0
    'AQ' = 0
    'AR' = 0
    'AS' = 0
    'AU' = 0
    'AV' = 0
    'AW' = 0
    'AX' = 0
    'AY' = 0
    'AZ' = 0
    'BA' = 0
    'BB' = 0
    'BC' = 0
    'BD' = 0
    'BE' = '0negawsklov0negawsklov'
    'BG' = 'HomeSweetHome'
    'BH' = 0
    'BI' = 0
}

# Per-row values that vary: A = reading datetime, AT = unix epoch, BF = wifiRssi.
$newRows = @(
    @{ Row = 14; A = '2022-03-18 12:38:47'; AT = 1647603529; BF = 54 }
    @{ Row = 15; A = '2022-03-18 12:39:39'; AT = 1647603581; BF = 50 }
    @{ Row = 16; A = '2022-03-18 12:54:16'; AT = 1647604458; BF = 52 }
    @{ Row = 17; A = '2022-03-18 12:58:39'; AT = 1647604721; BF = 48 }
    @{ Row = 18; A = '2022-03-18 12:59:52'; AT = 1647604794; BF = 50 }
    @{ Row = 19; A = '2022-03-18 13:02:04'; AT = 1647604926; BF = 44 }
    @{ Row = 20; A = '2022-03-18 13:07:37'; AT = 1647605259; BF = 48 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    foreach ($col in $rowTemplate.Keys) {
        $ws.Range("$col$r").Value = $rowTemplate[$col]
    }
    $ws.Range("A$r").Value = $entry.A
    $ws.Range("AT$r").Value = $entry.AT
    $ws.Range("BF$r").Value = $entry.BF
}
